$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: program_id changed from 2 to 3
$ws.Range("B4").Value = 3

# New row 5: a new segment for program_id 3
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "00-07-00_00-15-00.mov"

# Leave the cursor where the author last left it
$ws.Range("B5").Select()
